# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps to reflect a fresh report
# generation run, per commit "Generate Report for Handback".
#
# The ff0848a9-2451-4554-ac25-850074700392.md row's de-de xliff was
# regenerated, so both the Overview "Latest HO Xliff Generate Date" and the
# de-de sheet's "Correspond Handoff Datetime" (which shared that same
# timestamp) move to the new handoff time. The zh-cn handoff/handback times
# and the de-de handback time for that row also advance to a fresh run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-27 00:45:41"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-27 00:45:36"
$wsZhCn.Range("K4").Value = "2016-08-27 00:45:55"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-27 00:45:41"
$wsDeDe.Range("K4").Value = "2016-08-27 00:46:05"
